$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 401.8965782137555
$ws.Range("C2").Value = 736.8907799324985
$ws.Range("D2").Value = 383.60313119269222
$ws.Range("E2").Value = 639.9664850499895

$ws.Range("B3").Value = 404.9306028098606
$ws.Range("C3").Value = 508.21454833938918
$ws.Range("D3").Value = 490.67665434144442
$ws.Range("E3").Value = 405.41359672133422

$ws.Range("B1:E3").Select()
